# Renumber the "TransactionNumber" column (column D) shared-string values
# on Sheet2 and Sheet3, per the commit's data refresh.
#
# Sheet2 rows 2..51 (skipping 7, 17, 20, 30 which have no TransactionNumber
# cell) get sequential values starting at 986225.
# Sheet3 rows 2..51 (same skip pattern) get sequential values starting at
# 476054.
#
# Rows 7, 17, 20 and 30 correspond to the "Failed" / "Business Exception"
# sample rows, which use a different column layout and have no column D
# entry at all.

$wb = $excel.ActiveWorkbook

$skipRows = @(7, 17, 20, 30)

function Set-TransactionNumbers($sheetName, $startValue) {
    $ws = $wb.Worksheets.Item($sheetName)
    $counter = $startValue
    for ($row = 2; $row -le 51; $row++) {
        if ($skipRows -contains $row) {
            continue
        }
        $cell = $ws.Cells.Item($row, 4)
        # Force the new value to stay a text (shared-string) cell, matching
        # the existing "715431"-style text entries, then drop back to the
        # default style so no visible formatting changes.
        $cell.NumberFormat = "@"
        $cell.Value = [string]$counter
        $cell.Style = "Normal"
        $counter = $counter + 1
    }
}

Set-TransactionNumbers "Sheet2" 986225
Set-TransactionNumbers "Sheet3" 476054
